$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RF values (column I) for rows 32 through 58 from the old value
# (49.47355555555556) to the new value (47.11385714285715), reflecting
# the "Update of 2025 data and RF changes" commit.
$ws.Range("I32:I58").Value = 47.11385714285715
